# Realestate Update resale numbers 2024-01-24 09:06
# Append a new data row (row 90) to the CityResaleNum sheet with the
# latest resale-number snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 90

# Text columns: prefix with a leading apostrophe so Excel stores the
# literal text instead of auto-converting date/time-like strings or
# zero-padded numbers (e.g. "2024-01-24", "03") into dates/numbers.
$ws.Range("A$row").Value = "'2024-01-24"
$ws.Range("B$row").Value = "'09:06:58"
$ws.Range("C$row").Value = "'Wednesday"
$ws.Range("D$row").Value = "'03"

# Numeric columns (city resale numbers).
$ws.Range("E$row").Value = 138443
$ws.Range("F$row").Value = 141387
$ws.Range("G$row").Value = 170552
$ws.Range("H$row").Value = 148877
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 122985
$ws.Range("K$row").Value = 223591
$ws.Range("L$row").Value = 256148
$ws.Range("M$row").Value = 185041
$ws.Range("N$row").Value = 110028
$ws.Range("O$row").Value = 41303
$ws.Range("P$row").Value = 30889
$ws.Range("Q$row").Value = 73329
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 41927
$ws.Range("T$row").Value = -1
